$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1; Col = 1; Text = '55÷2=27, 1' },
    @{ Row = 1; Col = 2; Text = '12÷5=2, 2' },
    @{ Row = 1; Col = 3; Text = '85÷3=28, 1' },
    @{ Row = 1; Col = 4; Text = '31÷7=4, 3' },
    @{ Row = 1; Col = 5; Text = '66÷4=16, 2' },
    @{ Row = 5; Col = 1; Text = '92÷2=46, 0' },
    @{ Row = 5; Col = 2; Text = '35÷9=3, 8' },
    @{ Row = 5; Col = 3; Text = '36÷7=5, 1' },
    @{ Row = 5; Col = 4; Text = '18÷5=3, 3' },
    @{ Row = 5; Col = 5; Text = '69÷2=34, 1' },
    @{ Row = 9; Col = 1; Text = '61÷8=7, 5' },
    @{ Row = 9; Col = 2; Text = '90÷7=12, 6' },
    @{ Row = 9; Col = 3; Text = '50÷2=25, 0' },
    @{ Row = 9; Col = 4; Text = '25÷4=6, 1' },
    @{ Row = 9; Col = 5; Text = '51÷6=8, 3' },
    @{ Row = 13; Col = 1; Text = '73÷3=24, 1' },
    @{ Row = 13; Col = 2; Text = '72÷5=14, 2' },
    @{ Row = 13; Col = 3; Text = '71÷5=14, 1' },
    @{ Row = 13; Col = 4; Text = '59÷9=6, 5' },
    @{ Row = 13; Col = 5; Text = '19÷2=9, 1' },
    @{ Row = 17; Col = 1; Text = '46÷3=15, 1' },
    @{ Row = 17; Col = 2; Text = '25÷5=5, 0' },
    @{ Row = 17; Col = 3; Text = '91÷9=10, 1' },
    @{ Row = 17; Col = 4; Text = '52÷8=6, 4' },
    @{ Row = 17; Col = 5; Text = '19÷5=3, 4' }
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    $cell.Range.Text = $c.Text
}

Write-Output "Applied $($changes.Count) cell updates."
